# Normalize the "Recorded By" (column G) author-list ordering on the
# session-analysis sheet: for a handful of known value combinations, the
# first two comma-separated names were swapped (e.g. the synced author list
# now lists "System"/"system" before the human email address). This only
# touches the exact values below; other "Recorded By" combinations
# (e.g. "backup@backdoor.com, System" or solo values like
# "dnasr281@gmail.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "backup@backdoor.com, system, System" = "system, backup@backdoor.com, System"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $current = $cell.Value2
    if ($current -ne $null -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
